# Generate Report for Handback
# Updates the Correspond Handoff/Handback DateTime stamps for the
# zh-cn and de-de report rows (row 2 in each language sheet).

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-15 03:55:55"
$wsZhCn.Range("H2").Value = "2016-03-15 03:56:40"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-15 03:56:04"
$wsDeDe.Range("H2").Value = "2016-03-15 03:56:53"
